# "Change polyfit to ployval"
# The placeholder text "19..58" in B2 (stored as a shared string) is
# replaced with the actual computed numeric result of the fit, 19.58,
# and the active selection is moved from B34 to D25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.58

$ws.Range("D25").Select()
